# Update "time_taken" (column F) timestamps on the "data" sheet.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$timestamps = @(
  "2021-10-05 14:19:37.453905",
  "2021-10-05 14:19:37.453913",
  "2021-10-05 14:19:37.453916",
  "2021-10-05 14:19:37.453919",
  "2021-10-05 14:19:37.453922",
  "2021-10-05 14:19:37.453924",
  "2021-10-05 14:19:37.453927",
  "2021-10-05 14:19:37.453929",
  "2021-10-05 14:19:37.453932",
  "2021-10-05 14:19:37.453935",
  "2021-10-05 14:19:37.453937",
  "2021-10-05 14:19:37.453940",
  "2021-10-05 14:19:37.453942",
  "2021-10-05 14:19:37.453945",
  "2021-10-05 14:19:37.453947",
  "2021-10-05 14:19:37.453950",
  "2021-10-05 14:19:37.453952",
  "2021-10-05 14:19:37.453955",
  "2021-10-05 14:19:37.453958",
  "2021-10-05 14:19:37.453960",
  "2021-10-05 14:19:37.453962",
  "2021-10-05 14:19:37.453965",
  "2021-10-05 14:19:37.453968",
  "2021-10-05 14:19:37.453970",
  "2021-10-05 14:19:37.453973",
  "2021-10-05 14:19:37.453976",
  "2021-10-05 14:19:37.453978",
  "2021-10-05 14:19:37.453981",
  "2021-10-05 14:19:37.453983",
  "2021-10-05 14:19:37.453986",
  "2021-10-05 14:19:37.453989",
  "2021-10-05 14:19:37.453991",
  "2021-10-05 14:19:37.453994",
  "2021-10-05 14:19:37.453997",
  "2021-10-05 14:19:37.453999",
  "2021-10-05 14:19:37.454002",
  "2021-10-05 14:19:37.454004",
  "2021-10-05 14:19:37.454007",
  "2021-10-05 14:19:37.454009",
  "2021-10-05 14:19:37.454012",
  "2021-10-05 14:19:37.454015",
  "2021-10-05 14:19:37.454017",
  "2021-10-05 14:19:37.454020",
  "2021-10-05 14:19:37.454022",
  "2021-10-05 14:19:37.454025",
  "2021-10-05 14:19:37.454028",
  "2021-10-05 14:19:37.454030",
  "2021-10-05 14:19:37.454033",
  "2021-10-05 14:19:37.454035",
  "2021-10-05 14:19:37.454038",
  "2021-10-05 14:19:37.454040",
  "2021-10-05 14:19:37.454043",
  "2021-10-05 14:19:37.454046"
)

$row = 2
foreach ($ts in $timestamps) {
    $dataSheet.Range("F$row").Value = $ts
    $row = $row + 1
}

# Add a new "metadata" worksheet positioned right after "data".
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the bold/bordered header style already used by data!B1 so no new
# cell style gets introduced into the workbook.
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$newSheet.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

# Header row.
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "ClinGen_Familial thoracic aortic aneurysm and aortic dissection"
$newSheet.Range("C2").Value = 210

# "0.10" must stay a text value (not become the number 0.1). Enter it as a
# formula producing a string, then collapse the formula down to its
# literal value via copy / paste-values so no text-format style sticks.
$newSheet.Range("D2").Formula = '="0.10"'
$newSheet.Range("D2").Copy()
$newSheet.Range("D2").PasteSpecial(-4163)  # xlPasteValues

$newSheet.Range("E2").Value = "2017-11-05T02:37:20.232365Z"
$newSheet.Range("F2").Value = "2021-10-05 14:19:37.450470"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/210/?format=json"

$excel.CutCopyMode = 0
